$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row whose column A starts with the text of the post to be removed
# ("「四行目」") and delete the entire row. Excel automatically shifts all
# subsequent rows up and keeps the rest of the data intact.
$ws.Rows.Item(102).EntireRow.Delete()
